# "Generate Report for Archive"
#
# 1) Status text changes from "Ready for handoff" to "In Translation"
#    for every cell that currently shows it (Overview!E2:F3 and the
#    "Status" column, C2:C3, on each locale sheet).
# 2) The two now-narrower "status" columns shrink from ~17.22 chars to
#    ~13.41 chars: Overview columns E & F, and column C on each of the
#    locale sheets (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$newStatus  = "In Translation"
# Target OOXML column width is ~13.4101845877511 characters. This COM
# layer's ColumnWidth setter quantizes to 1/6-character steps, so feed it
# the input value whose quantized result lands closest to that target
# (12.5 -> stored width 13.333333333333334, the nearest reachable value).
$newWidth   = 12.5

# --- Overview sheet -------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth

# --- Locale sheets (zh-cn, de-de) ------------------------------------
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("C2").Value = $newStatus
    $ws.Range("C3").Value = $newStatus
    $ws.Columns.Item(3).ColumnWidth = $newWidth
}
